# Importar Datos Relacionados con Fast Excel y otro
# Adds an "acronym" column (E) to the users sheet, populated with a
# repeating area_a / area_b / area_c pattern for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# New header for column E
$ws.Range("E1").Value = "acronym"

# Fill column E (rows 2-11) with the repeating area_a/area_b/area_c values
$areas = @("area_a", "area_b", "area_c")
for ($i = 2; $i -le 11; $i++) {
    $area = $areas[($i - 2) % 3]
    $ws.Range("E$i").Value = $area
}

# New data cells wrap their text (matches the new cell style used for column E)
$ws.Range("E2:E11").WrapText = $true

# Leave the selection where the editor last left it
$ws.Range("A6").Select()
